$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__s = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.359.69'
$ws.Range('D2').Style = $__s
$__s = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.793.52'
$ws.Range('D3').Style = $__s
$__s = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E3').Style = $__s
$__s = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = $__s
$__s = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E4').Style = $__s
$__s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.007'
$ws.Range('D5').Style = $__s
$__s = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('E5').Style = $__s
$__s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.73'
$ws.Range('D6').Style = $__s
$__s = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('E6').Style = $__s
$__s = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('E7').Style = $__s
$__s = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3630'
$ws.Range('D8').Style = $__s
$__s = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('E8').Style = $__s
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$__s = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07071'
$ws.Range('D9').Style = $__s
$__s = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('E9').Style = $__s
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$__s = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8714'
$ws.Range('D10').Style = $__s
$__s = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E10').Style = $__s
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$__s = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07793'
$ws.Range('D11').Style = $__s
$__s = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('E11').Style = $__s
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$__s = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.51'
$ws.Range('D12').Style = $__s
$__s = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('E12').Style = $__s
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$__s = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.798.37'
$ws.Range('D13').Style = $__s
$__s = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.43%  '
$ws.Range('E13').Style = $__s
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$__s = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.270'
$ws.Range('D14').Style = $__s
$__s = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('E14').Style = $__s
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$__s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.311'
$ws.Range('D15').Style = $__s
$__s = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('E15').Style = $__s
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$__s = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '84.70'
$ws.Range('D16').Style = $__s
$__s = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.59%  '
$ws.Range('E16').Style = $__s
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$__s = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('D17').Style = $__s
$__s = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('E17').Style = $__s
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$__s = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008519'
$ws.Range('D18').Style = $__s
$__s = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('E18').Style = $__s
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$__s = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').Style = $__s
$__s = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('E19').Style = $__s
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$__s = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.391.90'
$ws.Range('D20').Style = $__s
$__s = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.07%  '
$ws.Range('E20').Style = $__s
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$__s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.17'
$ws.Range('D21').Style = $__s
$__s = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('E21').Style = $__s
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$__s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.977'
$ws.Range('D22').Style = $__s
$__s = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E22').Style = $__s
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$__s = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.47'
$ws.Range('D23').Style = $__s
$__s = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('E23').Style = $__s
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$__s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.975.28'
$ws.Range('D24').Style = $__s
$__s = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.87%  '
$ws.Range('E24').Style = $__s
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$__s = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.978'
$ws.Range('D25').Style = $__s
$__s = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E25').Style = $__s
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$__s = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.94'
$ws.Range('D26').Style = $__s
$__s = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('E26').Style = $__s
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$__s = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.86'
$ws.Range('D27').Style = $__s
$__s = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('E27').Style = $__s
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$__s = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.035'
$ws.Range('D28').Style = $__s
$__s = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E28').Style = $__s
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$__s = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.21'
$ws.Range('D29').Style = $__s
$__s = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('E29').Style = $__s
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$__s = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.827'
$ws.Range('D30').Style = $__s
$__s = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('E30').Style = $__s
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$__s = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08654'
$ws.Range('D31').Style = $__s
$__s = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.86%  '
$ws.Range('E31').Style = $__s
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$__s = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.025'
$ws.Range('D32').Style = $__s
$__s = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('E32').Style = $__s
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$__s = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.438'
$ws.Range('D33').Style = $__s
$__s = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('E33').Style = $__s
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$__s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7180'
$ws.Range('D34').Style = $__s
$__s = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.85%  '
$ws.Range('E34').Style = $__s
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$__s = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.111'
$ws.Range('D35').Style = $__s
$__s = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.27%  '
$ws.Range('E35').Style = $__s
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$__s = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.642'
$ws.Range('D36').Style = $__s
$__s = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('E36').Style = $__s
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$__s = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.006'
$ws.Range('D37').Style = $__s
$__s = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E37').Style = $__s
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$__s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.079'
$ws.Range('D38').Style = $__s
$__s = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E38').Style = $__s
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$__s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01940'
$ws.Range('D39').Style = $__s
$__s = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E39').Style = $__s
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$__s = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05089'
$ws.Range('D40').Style = $__s
$__s = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('E40').Style = $__s
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$__s = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.868'
$ws.Range('D41').Style = $__s
$__s = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('E41').Style = $__s
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$__s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5100'
$ws.Range('D42').Style = $__s
$__s = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('E42').Style = $__s
$__s = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.891'
$ws.Range('D43').Style = $__s
$__s = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('E43').Style = $__s
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$__s = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1518'
$ws.Range('D44').Style = $__s
$__s = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.09%  '
$ws.Range('E44').Style = $__s
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$__s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.976'
$ws.Range('D45').Style = $__s
$__s = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.50%  '
$ws.Range('E45').Style = $__s
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$__s = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.008'
$ws.Range('D46').Style = $__s
$__s = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E46').Style = $__s
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$__s = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4666'
$ws.Range('D47').Style = $__s
$__s = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('E47').Style = $__s
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$__s = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.845'
$ws.Range('D48').Style = $__s
$__s = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('E48').Style = $__s
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$__s = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '99.61'
$ws.Range('D49').Style = $__s
$__s = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.06%  '
$ws.Range('E49').Style = $__s
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$__s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.585'
$ws.Range('D50').Style = $__s
$__s = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('E50').Style = $__s
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$__s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05962'
$ws.Range('D51').Style = $__s
$__s = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.55%  '
$ws.Range('E51').Style = $__s
